$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 1.524740333333333
$ws.Range("H2").Value = 4.574221
$ws.Range("I2").Value = 0.2062237893390968
$ws.Range("J2").Value = 0.2062237893390969
$ws.Range("Q2").Value = 0.02030344227866666
$ws.Range("R2").Value = 0.182730980508
$ws.Range("S2").Value = 0.2062237893390968
$ws.Range("T2").Value = 0.2062237893390969

# Row 3
$ws.Range("I3").Value = 0.4308548451232278
$ws.Range("J3").Value = 0.4308548451232279
$ws.Range("R3").Value = 0.3817722899879999
$ws.Range("S3").Value = 0.4308548451232278
$ws.Range("T3").Value = 0.4308548451232279

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2032796666666667
$ws.Range("H4").Value = 0.609839
$ws.Range("I4").Value = 0.02749392945088694
$ws.Range("J4").Value = 0.02749392945088694
$ws.Range("Q4").Value = 0.002706872041333333
$ws.Range("R4").Value = 0.024361848372
$ws.Range("S4").Value = 0.02749392945088694
$ws.Range("T4").Value = 0.02749392945088694

# Row 5
$ws.Range("G5").Value = 2.061212666666667
$ws.Range("H5").Value = 6.183638
$ws.Range("I5").Value = 0.2787826080683977
$ws.Range("J5").Value = 0.2787826080683978
$ws.Range("Q5").Value = 0.02744710786933333
$ws.Range("R5").Value = 0.247023970824
$ws.Range("S5").Value = 0.2787826080683977
$ws.Range("T5").Value = 0.2787826080683978

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4188103333333333
$ws.Range("H6").Value = 1.256431
$ws.Range("I6").Value = 0.05664482801839063
$ws.Range("J6").Value = 0.05664482801839064
$ws.Range("Q6").Value = 0.005576878398666667
$ws.Range("R6").Value = 0.050191905588
$ws.Range("S6").Value = 0.05664482801839063
$ws.Range("T6").Value = 0.05664482801839064

$wb.Save()
